# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# D (Price) and E (Volume(1h)) columns are stored as TEXT in the source data
# (prices use "." as a thousands separator, e.g. "27.584.43", and volumes are
# padded percentage strings, e.g. "  -1.09%  "), so every write below is done in
# a way that keeps the cell type as text instead of letting Excel auto-convert
# number-looking price strings (e.g. "1.99") into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.584.43"
$ws.Range("E2").Value = "  -1.09%  "

# Row 3
$ws.Range("D3").Value = "1.596.51"
$ws.Range("E3").Value = "  -2.10%  "

# Row 4
$ws.Range("E4").Value = "  +0.36%  "

# Row 5
$ws.Range("D5").Formula = "=""207.97"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E5").Value = "  -1.62%  "

# Row 6
$ws.Range("E6").Value = "  -3.79%  "

# Row 7
$ws.Range("E7").Value = "  +0.40%  "

# Row 8
$ws.Range("D8").Formula = "=""22.35"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E8").Value = "  -4.29%  "

# Row 9
$ws.Range("E9").Value = "  -1.98%  "

# Row 10
$ws.Range("E10").Value = "  -3.27%  "

# Row 11
$ws.Range("E11").Value = "  -1.87%  "

# Row 12
$ws.Range("D12").Value = "1.823.92"
$ws.Range("E12").Value = "  -2.09%  "

# Row 13
$ws.Range("D13").Value = "1.621.23"
$ws.Range("E13").Value = "  -0.61%  "

# Row 14
$ws.Range("E14").Value = "  -4.07%  "

# Row 15
$ws.Range("E15").Value = "  -4.40%  "

# Row 16
$ws.Range("D16").Formula = "=""63.41"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E16").Value = "  -2.96%  "

# Row 17
$ws.Range("D17").Value = "27.599.59"
$ws.Range("E17").Value = "  -1.04%  "

# Row 18
$ws.Range("D18").Formula = "=""217.44"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E18").Value = "  -5.21%  "

# Row 19
$ws.Range("E19").Value = "  -4.00%  "

# Row 20
$ws.Range("E20").Value = "  -3.55%  "

# Row 21
$ws.Range("E21").Value = "  +0.42%  "

# Row 22
$ws.Range("E22").Value = "  -3.77%  "

# Row 23
$ws.Range("D23").Formula = "=""9.61"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E23").Value = "  -4.57%  "

# Row 24
$ws.Range("D24").Formula = "=""1.99"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E24").Value = "  -3.35%  "

# Row 25
$ws.Range("D25").Formula = "=""152.62"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E25").Value = "  -1.09%  "

# Row 26
$ws.Range("E26").Value = "  -1.73%  "

# Row 27
$ws.Range("E27").Value = "  +0.42%  "

# Row 28
$ws.Range("D28").Formula = "=""15.09"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E28").Value = "  -2.69%  "

# Row 30
$ws.Range("E30").Value = "  -1.79%  "

# Row 31
$ws.Range("E31").Value = "  -2.97%  "

# Row 32
$ws.Range("D32").Formula = "=""3.27"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E32").Value = "  -4.49%  "

# Row 33
$ws.Range("D33").Value = "1.375.08"
$ws.Range("E33").Value = "  -0.98%  "

# Row 34
$ws.Range("D34").Formula = "=""2.95"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E34").Value = "  -5.13%  "

# Row 35
$ws.Range("E35").Value = "  -3.96%  "

# Row 36
$ws.Range("D36").Formula = "=""0.967"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E36").Value = "  -4.81%  "

# Row 37
$ws.Range("E37").Value = "  -1.25%  "

# Row 38
$ws.Range("D38").Formula = "=""0.0165"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E38").Value = "  -3.39%  "

# Row 39
$ws.Range("D39").Formula = "=""0.541"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E39").Value = "  -3.01%  "

# Row 40
$ws.Range("D40").Formula = "=""0.812"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E40").Value = "  -4.41%  "

# Row 41
$ws.Range("E41").Value = "  +0.38%  "

# Row 42
$ws.Range("E42").Value = "  -3.88%  "

# Row 43
$ws.Range("D43").Formula = "=""5.36"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E43").Value = "  -1.10%  "

# Row 44
$ws.Range("D44").Formula = "=""1.78"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E44").Value = "  -2.95%  "

# Row 45
$ws.Range("D45").Formula = "=""64.06"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E45").Value = "  -2.39%  "

# Row 46
$ws.Range("D46").Formula = "=""2.18"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E46").Value = "  +1.39%  "

# Row 47
$ws.Range("D47").Value = "1.733.19"
$ws.Range("E47").Value = "  -2.22%  "

# Row 48
$ws.Range("D48").Formula = "=""87.21"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("E48").Value = "  -1.30%  "

# Row 49
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  -3.01%  "

# Row 50
$ws.Range("E50").Value = "  -4.35%  "

# Row 51
$ws.Range("E51").Value = "  -1.20%  "
